$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = " Zinc"
$ws.Range("A3").Value = "Copper"
$ws.Range("A4").Value = "Diesel"

$ws.Range("A5").Value = "Steel Products"
$ws.Range("B5").Value = "1"

$ws.Range("A6").Value = "Titanium Dioxide"
$ws.Range("B6").Value = "1"

$ws.Range("A7").Value = "Steel — Hot Rolled"
$ws.Range("B7").Value = "1"

$ws.Range("A8").Value = "Precious Metals"
$ws.Range("B8").Value = "1"

$ws.Range("A9").Value = "Solvents"
$ws.Range("B9").Value = "1"

$ws.Range("A10").Value = "Steel — Carbon"
$ws.Range("B10").Value = "1"

$ws.Range("B11").Value = "10"

$ws.Range("A12").Value = "Corrugated Boxes"
$ws.Range("B12").Value = "2"

$ws.Range("A13").Value = "Corrugated Sheets"
$ws.Range("B13").Value = "2"

$ws.Range("A14").Value = "Crude Oil"
$ws.Range("B14").Value = "2"

$ws.Range("A15").Value = "Gasoline"
$ws.Range("B15").Value = "2"

$ws.Range("A16").Value = "Plastic Resins"
$ws.Range("B16").Value = "4"

$ws.Range("A17").Value = "Aluminum"
$ws.Range("B17").Value = "5"

$ws.Range("A18").Value = "Polypropylene"
$ws.Range("B18").Value = "7"

$ws.Range("A19").Value = "High-Density Polyethylene"
$ws.Range("B19").Value = "HDPE) Resin"
